$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 110 is a brand new row appended at the end of the table, holding the
# data that previously lived in row 109. Populate the columns that are
# constant across every data row in this sheet (A, B, C, E, F, G, H, I, N,
# O, Q, R) with the same literal values used throughout the table.
$ws.Cells.Item(110, 1).Value = 8
$ws.Cells.Item(110, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(110, 3).Value = "Coquimbo"
$ws.Cells.Item(110, 5).Value = 4
$ws.Cells.Item(110, 6).Value = 100112040
$ws.Cells.Item(110, 7).Value = "Cilantro"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(110, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(110, 17).Value = 1.5
$ws.Cells.Item(110, 18).Value = "Hortaliza"
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(109, 4).NumberFormat

# Data for rows 19-110: the weekly price table shifted down by one row to
# make room for a new, more recent reading in row 19, pushing the former
# last row (old row 109) down into the newly appended row 110.
$rows = @(
  @{Row=19; D=44560; J=3600; K=2500; L=3000; M=2750; P=1833},
  @{Row=20; D=44159; J=2900; K=1000; L=1500; M=1250; P=833},
  @{Row=21; D=44336; J=3600; K=1300; L=1500; M=1400; P=933},
  @{Row=22; D=44343; J=3600; K=1300; L=1500; M=1400; P=933},
  @{Row=23; D=44245; J=3500; K=1500; L=2000; M=1750; P=1167},
  @{Row=24; D=44397; J=3300; K=1500; L=2000; M=1750; P=1167},
  @{Row=25; D=44523; J=3280; K=1500; L=2000; M=1750; P=1167},
  @{Row=26; D=44369; J=3360; K=1500; L=2000; M=1750; P=1167},
  @{Row=27; D=44189; J=3200; K=1400; L=1500; M=1450; P=967},
  @{Row=28; D=44355; J=3400; K=1300; L=1500; M=1400; P=933},
  @{Row=29; D=44406; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=30; D=44215; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=31; D=44392; J=3540; K=1500; L=2000; M=1750; P=1167},
  @{Row=32; D=44544; J=3200; K=1500; L=2000; M=1750; P=1167},
  @{Row=33; D=44201; J=2800; K=1300; L=1500; M=1400; P=933},
  @{Row=34; D=44530; J=3300; K=1500; L=2000; M=1750; P=1167},
  @{Row=35; D=44203; J=3300; K=1300; L=1500; M=1400; P=933},
  @{Row=36; D=44236; J=3400; K=1500; L=2000; M=1750; P=1167},
  @{Row=37; D=44537; J=3300; K=1500; L=2000; M=1750; P=1167},
  @{Row=38; D=44320; J=3600; K=1300; L=1500; M=1400; P=933},
  @{Row=39; D=44294; J=3500; K=2000; L=2500; M=2250; P=1500},
  @{Row=40; D=44315; J=3600; K=1300; L=1500; M=1400; P=933},
  @{Row=41; D=44483; J=3160; K=1500; L=2000; M=1750; P=1167},
  @{Row=42; D=44348; J=3560; K=1300; L=1500; M=1400; P=933},
  @{Row=43; D=44488; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=44; D=44546; J=2800; K=2000; L=2500; M=2250; P=1500},
  @{Row=45; D=44252; J=3400; K=1500; L=2000; M=1750; P=1167},
  @{Row=46; D=44383; J=3300; K=1500; L=2000; M=1750; P=1167},
  @{Row=47; D=44222; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=48; D=44376; J=3480; K=1500; L=2000; M=1750; P=1167},
  @{Row=49; D=44243; J=3400; K=1500; L=2000; M=1750; P=1167},
  @{Row=50; D=44292; J=3500; K=2000; L=2500; M=2250; P=1500},
  @{Row=51; D=44299; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=52; D=44166; J=2900; K=1300; L=1500; M=1400; P=933},
  @{Row=53; D=44327; J=3540; K=1300; L=1500; M=1400; P=933},
  @{Row=54; D=44390; J=3320; K=1500; L=2000; M=1750; P=1167},
  @{Row=55; D=44273; J=3400; K=2000; L=2500; M=2250; P=1500},
  @{Row=56; D=44364; J=3600; K=1500; L=2000; M=1750; P=1167},
  @{Row=57; D=44469; J=3160; K=1500; L=2000; M=1750; P=1167},
  @{Row=58; D=44280; J=3400; K=2000; L=2500; M=2250; P=1500},
  @{Row=59; D=44525; J=3100; K=1500; L=2000; M=1750; P=1167},
  @{Row=60; D=44266; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=61; D=44306; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=62; D=44516; J=3300; K=1300; L=1500; M=1400; P=933},
  @{Row=63; D=44427; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=64; D=44495; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=65; D=44539; J=3200; K=1500; L=2000; M=1750; P=1167},
  @{Row=66; D=44532; J=3200; K=1800; L=2000; M=1900; P=1267},
  @{Row=67; D=44210; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=68; D=44271; J=3200; K=2000; L=2500; M=2250; P=1500},
  @{Row=69; D=44425; J=3360; K=2000; L=2500; M=2250; P=1500},
  @{Row=70; D=44476; J=3080; K=1500; L=2000; M=1750; P=1167},
  @{Row=71; D=44250; J=3600; K=1500; L=2000; M=1750; P=1167},
  @{Row=72; D=44168; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=73; D=44553; J=3300; K=2000; L=2500; M=2250; P=1500},
  @{Row=74; D=44161; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=75; D=44285; J=3400; K=2000; L=2500; M=2250; P=1500},
  @{Row=76; D=44418; J=3400; K=2000; L=2500; M=2250; P=1500},
  @{Row=77; D=44434; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=78; D=44467; J=3200; K=1500; L=2000; M=1750; P=1167},
  @{Row=79; D=44231; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=80; D=44490; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=81; D=44259; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=82; D=44341; J=3560; K=1300; L=1500; M=1400; P=933},
  @{Row=83; D=44208; J=3200; K=1400; L=1500; M=1450; P=967},
  @{Row=84; D=44264; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=85; D=44322; J=3640; K=1300; L=1500; M=1400; P=933},
  @{Row=86; D=44551; J=3200; K=2000; L=2500; M=2250; P=1500},
  @{Row=87; D=44420; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=88; D=44385; J=3560; K=1500; L=2000; M=1750; P=1167},
  @{Row=89; D=44278; J=3200; K=2000; L=2500; M=2250; P=1500},
  @{Row=90; D=44308; J=3600; K=1300; L=1500; M=1400; P=933},
  @{Row=91; D=44187; J=3200; K=1400; L=1500; M=1450; P=967},
  @{Row=92; D=44474; J=3200; K=1500; L=2000; M=1750; P=1167},
  @{Row=93; D=44446; J=3260; K=2000; L=2500; M=2250; P=1500},
  @{Row=94; D=44350; J=3600; K=1200; L=1500; M=1350; P=900},
  @{Row=95; D=44413; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=96; D=44238; J=3400; K=1500; L=2000; M=1750; P=1167},
  @{Row=97; D=44257; J=3600; K=2000; L=2500; M=2250; P=1500},
  @{Row=98; D=44411; J=3400; K=2000; L=2500; M=2250; P=1500},
  @{Row=99; D=44175; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=100; D=44196; J=3200; K=1400; L=1500; M=1450; P=967},
  @{Row=101; D=44432; J=3360; K=2000; L=2500; M=2250; P=1500},
  @{Row=102; D=44224; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=103; D=44329; J=3560; K=1300; L=1500; M=1400; P=933},
  @{Row=104; D=44511; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=105; D=44357; J=3600; K=1300; L=1500; M=1400; P=933},
  @{Row=106; D=44371; J=3520; K=1500; L=2000; M=1750; P=1167},
  @{Row=107; D=44194; J=3200; K=1400; L=1500; M=1450; P=967},
  @{Row=108; D=44313; J=3700; K=1300; L=1500; M=1400; P=933},
  @{Row=109; D=44518; J=3200; K=1300; L=1500; M=1400; P=933},
  @{Row=110; D=44540; J=3000; K=1500; L=2000; M=1750; P=1167}
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 10).Value = $item.J
    $ws.Cells.Item($item.Row, 11).Value = $item.K
    $ws.Cells.Item($item.Row, 12).Value = $item.L
    $ws.Cells.Item($item.Row, 13).Value = $item.M
    $ws.Cells.Item($item.Row, 16).Value = $item.P
}

Write-Host "Update complete"
